# waveform4R.xlsx update — add new "added date" column/header, two new
# parameter rows (t_1sE2.exp / t_1E2.exp), and clear a stale note.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# --- New header for column M ("Cited in date added" column) ---------------
$ws.Range("M1").Value = "f_added_date"

# --- Clear the stale "filter(max(ff89))" function note on row 20 ----------
$ws.Range("L20").ClearContents()

# --- New row 30: t_1sE2.exp -------------------------------------------------
$ws.Range("A30").Value = 109
$ws.Range("B30").Value = "t_1sE2.exp"
$ws.Range("C30").Value = "Sequential"
$ws.Range("D30").Value = "Time from 1st probe to 1st sustained E2 (10 minutes)" + [char]0x2020
$ws.Range("E30").Value = "seconds"
$ws.Range("F30").Value = "All tissues"
$ws.Range("G30").Value = "Plant acceptability"
$ws.Range("H30").Value = "value"
$ws.Range("I30").Value = "4, 5"
$ws.Range("J30").Value = "time to first sustained E2 - first 2 wave"
$ws.Range("L30").Value = "f109"
$ws.Range("M30").Value = (Get-Date -Year 2020 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M30").NumberFormat = "m/d/yyyy"

# --- New row 31: t_1E2.exp --------------------------------------------------
$ws.Range("A31").Value = 112
$ws.Range("B31").Value = "t_1E2.exp"
$ws.Range("C31").Value = "Sequential"
$ws.Range("D31").Value = "Time from 1st probe to 1st E2" + [char]0x2020
$ws.Range("E31").Value = "seconds"
$ws.Range("F31").Value = "All tissues"
$ws.Range("G31").Value = "Plant acceptability"
$ws.Range("H31").Value = "value"
$ws.Range("I31").Value = "4, 5 "
$ws.Range("J31").Value = "time to first E2 - first 2 wave"
$ws.Range("L31").Value = "f112"
$ws.Range("M31").Value = (Get-Date -Year 2020 -Month 9 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M31").NumberFormat = "m/d/yyyy"

# --- Match formatting of the rest of the table for the two new rows -------
$ws.Range("A30:M31").WrapText = $true
$ws.Range("A30:M31").VerticalAlignment = -4108  # xlCenter
$ws.Range("A30:M31").HorizontalAlignment = -4131  # xlLeft

$ws.Range("A30:M31").RowHeight = 60

Write-Host "waveform4R.xlsx update applied"
